$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay as Text so numeric-looking strings (e.g. "9.10") are not
# silently coerced into Number values by Excel's type inference on assignment.
foreach ($ref in @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "E9", "E10", "E11", "E12", "D13", "E13", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "E24", "E25", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "D31", "E31", "E32", "E33", "E34", "D35", "E37", "E38", "E39", "E40", "E41", "D42", "E42", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "D49", "E49", "B50", "C50", "D50", "E50", "B51", "C51", "D51", "E51")) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.056.25'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '2.298.38'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '300.38'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = '98.14'
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("E7").Value = '  +2.00%  '
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").Value = '17.71'
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").Value = '2.655.97'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").Value = '2.290.86'
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("D17").Value = '0.787'
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("D18").Value = '42.931.47'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = '12.77'
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = '0.0₃0914'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").Value = '69.01'
$ws.Range("E22").Value = '  +1.49%  '
$ws.Range("D23").Value = '237.29'
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E24").Value = '  -3.86%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -1.11%  '
$ws.Range("D27").Value = '24.97'
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("D28").Value = '165.25'
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '9.10'
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").Value = '33.13'
$ws.Range("E31").Value = '  -4.09%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  +2.21%  '
$ws.Range("D35").Value = '17.94'
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("E41").Value = '  +0.36%  '
$ws.Range("D42").Value = '2.011.84'
$ws.Range("E42").Value = '  +1.23%  '
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").Value = '10.28'
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("D46").Value = '17.44'
$ws.Range("E46").Value = '  -1.23%  '
$ws.Range("D47").Value = '2.82'
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("D49").Value = '2.522.86'
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = '73.32'
$ws.Range("E50").Value = '  +3.49%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '1.53'
$ws.Range("E51").Value = '  -1.89%  '
